$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range("D2").Value = "63.887.27"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "3.117.23"
$ws.Range("E3").Value = "  +0.12%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "585.87"
$ws.Range("E5").Value = "  -0.26%  "
Set-TextValue "D6" "146.41"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.112.27"
$ws.Range("E8").Value = "  +0.31%  "
Set-TextValue "D9" "0.531"
$ws.Range("E9").Value = "  -0.20%  "
Set-TextValue "D10" "0.160"
$ws.Range("E10").Value = "  +9.53%  "
Set-TextValue "D11" "5.77"
$ws.Range("E11").Value = "  +1.50%  "
Set-TextValue "D12" "0.463"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("E13").Value = "  +2.55%  "
Set-TextValue "D14" "37.24"
$ws.Range("E14").Value = "  +4.69%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "3.632.76"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "63.751.70"
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "3.121.92"
$ws.Range("E19").Value = "  +0.34%  "
Set-TextValue "D20" "464.70"
$ws.Range("E20").Value = "  +2.12%  "
Set-TextValue "D21" "14.34"
$ws.Range("E21").Value = "  +1.76%  "
$ws.Range("E22").Value = "  -0.52%  "
Set-TextValue "D23" "7.53"
$ws.Range("E23").Value = "  -1.05%  "
Set-TextValue "D24" "13.18"
$ws.Range("E24").Value = "  -3.49%  "
Set-TextValue "D25" "81.84"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -0.13%  "
Set-TextValue "D27" "8.93"
$ws.Range("E27").Value = "  +7.46%  "
Set-TextValue "D28" "2.70"
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  -0.33%  "
Set-TextValue "D32" "27.01"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("E33").Value = "  -2.61%  "
$ws.Range("D34").Value = "0.0₃0868"
$ws.Range("E34").Value = "  +5.58%  "
$ws.Range("B35").Value = "Stacks"
$ws.Range("C35").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D35" "2.36"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D36" "1.05"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +10.12%  "
$ws.Range("E38").Value = "  -0.02%  "
Set-TextValue "D39" "50.98"
$ws.Range("E39").Value = "  +0.01%  "
Set-TextValue "D40" "448.19"
$ws.Range("E40").Value = "  +4.57%  "
Set-TextValue "D41" "8.69"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "2.878.22"
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  +0.63%  "
$ws.Range("E45").Value = "  -0.83%  "
Set-TextValue "D46" "2.17"
$ws.Range("E46").Value = "  -0.17%  "
Set-TextValue "D47" "35.79"
$ws.Range("E47").Value = "  +3.08%  "
Set-TextValue "D49" "123.70"
$ws.Range("E49").Value = "  -1.06%  "
$ws.Range("E50").Value = "  -0.53%  "
Set-TextValue "D51" "24.68"
$ws.Range("E51").Value = "  -0.82%  "
